$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.155790090560913
$ws.Range("B1").Value = 2.383633136749268
$ws.Range("D1").Value = 2.391962289810181
$ws.Range("E1").Value = 1.224308371543884
